$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.197.19'
$ws.Range("E2").Value = '  -2.79%  '
$ws.Range("D3").Value = '3.012.67'
$ws.Range("E3").Value = '  -1.79%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.76'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.00'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.009.12'
$ws.Range("E8").Value = '  -1.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.499'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.29%  '
$ws.Range("E10").Value = '  -3.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.07'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.447'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000221'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.39'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.80%  '
$ws.Range("D15").Value = '3.503.15'
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").Value = '61.250.15'
$ws.Range("E17").Value = '  -2.68%  '
$ws.Range("D18").Value = '3.013.18'
$ws.Range("E18").Value = '  -1.82%  '
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '465.85'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.26'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.677'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.96'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.49'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.07'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.68'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.84'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("E30").Value = '  +0.90%  '
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("E32").Value = '  +2.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.49'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '55.57'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.29'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.91'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '458.99'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.87%  '
$ws.Range("D38").Value = '3.222.41'
$ws.Range("E38").Value = '  +4.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0787'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("E40").Value = '  -2.06%  '
$ws.Range("E41").Value = '  +2.85%  '
$ws.Range("E42").Value = '  +1.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '27.61'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +13.42%  '
$ws.Range("E44").Value = '  -4.65%  '
$ws.Range("E46").Value = '  -2.05%  '
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.96'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.96%  '
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("D50").Value = '0.0₃0495'
$ws.Range("E50").Value = '  -7.98%  '
$ws.Range("E51").Value = '  +7.84%  '
